# Apply edit: reorder model names in column A (rows 2..26) and set all
# metric columns (B..I) to the values of the best model (model_11_9_0).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @(
    "model_11_9_0",
    "model_11_9_22",
    "model_11_9_21",
    "model_11_9_20",
    "model_11_9_19",
    "model_11_9_18",
    "model_11_9_17",
    "model_11_9_16",
    "model_11_9_15",
    "model_11_9_14",
    "model_11_9_13",
    "model_11_9_23",
    "model_11_9_12",
    "model_11_9_10",
    "model_11_9_9",
    "model_11_9_8",
    "model_11_9_7",
    "model_11_9_6",
    "model_11_9_5",
    "model_11_9_4",
    "model_11_9_3",
    "model_11_9_2",
    "model_11_9_1",
    "model_11_9_11",
    "model_11_9_24"
)

$vals = @(
    0.3494677884409869,
    -0.7266016666444675,
    0.2087795153652025,
    -0.08449119809865491,
    0.7199474573135376,
    0.9218135476112366,
    0.9635273814201355,
    0.9414435029029846
)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]
    for ($c = 0; $c -lt $vals.Length; $c++) {
        $ws.Cells.Item($row, $c + 2).Value = $vals[$c]
    }
}
